$d = $word.ActiveDocument

$d.Content.Find.Execute("357×7=", $true, $false, $false, $false, $false, $true, 1, $false, "221×4=", 2) | Out-Null
$d.Content.Find.Execute("625×8=", $true, $false, $false, $false, $false, $true, 1, $false, "772×3=", 2) | Out-Null
$d.Content.Find.Execute("106×2=", $true, $false, $false, $false, $false, $true, 1, $false, "333×4=", 2) | Out-Null
$d.Content.Find.Execute("522×2=", $true, $false, $false, $false, $false, $true, 1, $false, "332×5=", 2) | Out-Null
$d.Content.Find.Execute("340×3=", $true, $false, $false, $false, $false, $true, 1, $false, "392×6=", 2) | Out-Null
$d.Content.Find.Execute("896×9=", $true, $false, $false, $false, $false, $true, 1, $false, "237×4=", 2) | Out-Null
$d.Content.Find.Execute("455×5=", $true, $false, $false, $false, $false, $true, 1, $false, "248×3=", 2) | Out-Null
$d.Content.Find.Execute("522×6=", $true, $false, $false, $false, $false, $true, 1, $false, "534×6=", 2) | Out-Null
$d.Content.Find.Execute("260×8=", $true, $false, $false, $false, $false, $true, 1, $false, "528×8=", 2) | Out-Null
$d.Content.Find.Execute("639×9=", $true, $false, $false, $false, $false, $true, 1, $false, "526×9=", 2) | Out-Null
$d.Content.Find.Execute("484×2=", $true, $false, $false, $false, $false, $true, 1, $false, "989×8=", 2) | Out-Null
$d.Content.Find.Execute("577×7=", $true, $false, $false, $false, $false, $true, 1, $false, "607×5=", 2) | Out-Null
$d.Content.Find.Execute("338×9=", $true, $false, $false, $false, $false, $true, 1, $false, "892×2=", 2) | Out-Null
$d.Content.Find.Execute("441×9=", $true, $false, $false, $false, $false, $true, 1, $false, "174×5=", 2) | Out-Null
$d.Content.Find.Execute("611×5=", $true, $false, $false, $false, $false, $true, 1, $false, "941×2=", 2) | Out-Null
$d.Content.Find.Execute("167×4=", $true, $false, $false, $false, $false, $true, 1, $false, "291×9=", 2) | Out-Null
$d.Content.Find.Execute("153×8=", $true, $false, $false, $false, $false, $true, 1, $false, "118×3=", 2) | Out-Null
$d.Content.Find.Execute("693×4=", $true, $false, $false, $false, $false, $true, 1, $false, "753×4=", 2) | Out-Null
$d.Content.Find.Execute("114×7=", $true, $false, $false, $false, $false, $true, 1, $false, "532×8=", 2) | Out-Null
$d.Content.Find.Execute("667×9=", $true, $false, $false, $false, $false, $true, 1, $false, "406×8=", 2) | Out-Null
$d.Content.Find.Execute("497×4=", $true, $false, $false, $false, $false, $true, 1, $false, "815×7=", 2) | Out-Null
$d.Content.Find.Execute("688×7=", $true, $false, $false, $false, $false, $true, 1, $false, "333×7=", 2) | Out-Null
$d.Content.Find.Execute("776×2=", $true, $false, $false, $false, $false, $true, 1, $false, "582×8=", 2) | Out-Null
$d.Content.Find.Execute("785×9=", $true, $false, $false, $false, $false, $true, 1, $false, "426×8=", 2) | Out-Null
$d.Content.Find.Execute("115×6=", $true, $false, $false, $false, $false, $true, 1, $false, "327×9=", 2) | Out-Null
